$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "26.142.96"
Set-TextValue "E2" "  -0.41%  "

# Row 3
Set-TextValue "D3" "1.650.10"
Set-TextValue "E3" "  -0.51%  "

# Row 4
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.19%  "

# Row 5
Set-TextValue "D5" "218.64"
Set-TextValue "E5" "  -0.26%  "

# Row 6
Set-TextValue "D6" "0.5212"
Set-TextValue "E6" "  -0.61%  "

# Row 7
Set-TextValue "E7" "  -0.23%  "

# Row 8
Set-TextValue "D8" "0.2623"
Set-TextValue "E8" "  -0.85%  "

# Row 9
Set-TextValue "D9" "0.06307"
Set-TextValue "E9" "  -0.23%  "

# Row 10
Set-TextValue "D10" "20.41"
Set-TextValue "E10" "  -1.38%  "

# Row 11
Set-TextValue "D11" "0.07649"
Set-TextValue "E11" "  -1.77%  "

# Row 12
Set-TextValue "D12" "4.593"
Set-TextValue "E12" "  +1.75%  "

# Row 13
Set-TextValue "D13" "1.647.29"
Set-TextValue "E13" "  +2.21%  "

# Row 14
Set-TextValue "D14" "1.874.25"
Set-TextValue "E14" "  -0.71%  "

# Row 15
Set-TextValue "D15" "0.5588"
Set-TextValue "E15" "  -0.82%  "

# Row 16
Set-TextValue "D16" "0.0₅8137"
Set-TextValue "E16" "  +0.83%  "

# Row 17
Set-TextValue "D17" "65.16"
Set-TextValue "E17" "  -0.24%  "

# Row 18
Set-TextValue "D18" "26.076.95"
Set-TextValue "E18" "  -0.65%  "

# Row 19
Set-TextValue "E19" "  -0.14%  "

# Row 20
Set-TextValue "D20" "4.598"
Set-TextValue "E20" "  -2.73%  "

# Row 21
Set-TextValue "D21" "194.40"
Set-TextValue "E21" "  +0.04%  "

# Row 22
Set-TextValue "D22" "10.47"
Set-TextValue "E22" "  +2.15%  "

# Row 23
Set-TextValue "D23" "5.929"
Set-TextValue "E23" "  -1.78%  "

# Row 24
Set-TextValue "D24" "1.003"
Set-TextValue "E24" "  -0.14%  "

# Row 25
Set-TextValue "D25" "145.19"
Set-TextValue "E25" "  -0.12%  "

# Row 26
Set-TextValue "D26" "0.1185"
Set-TextValue "E26" "  -2.18%  "

# Row 27
Set-TextValue "D27" "7.204"
Set-TextValue "E27" "  -0.39%  "

# Row 28
Set-TextValue "B28" "EthereumClassic"
Set-TextValue "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "15.86"
Set-TextValue "E28" "  -1.12%  "

# Row 29
Set-TextValue "B29" "Toncoin"
Set-TextValue "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "1.530"
Set-TextValue "E29" "  +1.72%  "

# Row 30
Set-TextValue "D30" "0.05465"
Set-TextValue "E30" "  -3.15%  "

# Row 31
Set-TextValue "E31" "  -0.76%  "

# Row 32
Set-TextValue "D32" "3.437"
Set-TextValue "E32" "  -1.48%  "

# Row 33
Set-TextValue "D33" "3.333"
Set-TextValue "E33" "  -1.26%  "

# Row 34
Set-TextValue "D34" "1.562"
Set-TextValue "E34" "  -2.60%  "

# Row 35
Set-TextValue "D35" "2.411"
Set-TextValue "E35" "  +0.30%  "

# Row 36
Set-TextValue "E36" "  -0.95%  "

# Row 37
Set-TextValue "D37" "0.9445"
Set-TextValue "E37" "  -0.22%  "

# Row 38
Set-TextValue "D38" "0.5609"
Set-TextValue "E38" "  -2.56%  "

# Row 39
Set-TextValue "D39" "0.01571"
Set-TextValue "E39" "  -2.17%  "

# Row 40
Set-TextValue "E40" "  -0.13%  "

# Row 41
Set-TextValue "D41" "5.737"
Set-TextValue "E41" "  -4.10%  "

# Row 42
Set-TextValue "D42" "1.027.79"
Set-TextValue "E42" "  -1.89%  "

# Row 43
Set-TextValue "D43" "0.8197"
Set-TextValue "E43" "  -3.19%  "

# Row 44
Set-TextValue "D44" "100.63"
Set-TextValue "E44" "  -2.18%  "

# Row 45
Set-TextValue "D45" "1.786.62"
Set-TextValue "E45" "  -0.69%  "

# Row 46
Set-TextValue "B46" "Aave"
Set-TextValue "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "57.26"
Set-TextValue "E46" "  -2.06%  "

# Row 47
Set-TextValue "B47" "BabyDogeCoin"
Set-TextValue "C47" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D47" "0.0₈106"
Set-TextValue "E47" "  +1.24%  "

# Row 48
Set-TextValue "D48" "1.000"
Set-TextValue "E48" "  -0.19%  "

# Row 49
Set-TextValue "D49" "0.4325"
Set-TextValue "E49" "  -0.67%  "

# Row 50
Set-TextValue "D50" "7.894"
Set-TextValue "E50" "  -1.68%  "

# Row 51
Set-TextValue "D51" "0.05116"
Set-TextValue "E51" "  -4.04%  "
